$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Apply the "mtitleStyle" (s=4, bold black Century w/ thin border) to the
#    three row-label cells in the summary block, which previously had the
#    default style. Donor: A9 already carries style index 4.
# ---------------------------------------------------------------------------
$ws.Range("A9").Copy()
$ws.Range("A10:A12").PasteSpecial(-4122)   # xlPasteFormats

# ---------------------------------------------------------------------------
# 2) Update the summary numbers (Right / Wrong / Not Attempt / Max / total).
# ---------------------------------------------------------------------------
$ws.Range("B10").Value = 3
$ws.Range("C10").Value = 1
$ws.Range("D10").Value = 24
$ws.Range("E10").Value = 28

$ws.Range("B11").Value = 4
$ws.Range("C11").Value = -1   # now a real number instead of a text "-1"

$ws.Range("B12").Value = 12
$ws.Range("C12").Value = -1
$ws.Range("E12").Value = "11/112"   # was "Absent"

# ---------------------------------------------------------------------------
# 3) Drop the third "Student Ans / Correct Ans" answer block entirely
#    (columns G:H, rows 15-40).
# ---------------------------------------------------------------------------
$ws.Range("G15:H40").Clear()

# ---------------------------------------------------------------------------
# 4) Consolidate the second answer block (columns D:E) into the first
#    (columns A:B) for a handful of rows, then drop the rest of column D:E.
# ---------------------------------------------------------------------------

# Row 17: D17 picks up "Option C" using the "normalStyle" (s=5) already used
# by the numeric "Right" column; E17 already holds "Option C".
$ws.Range("B10").Copy()
$ws.Range("D17").PasteSpecial(-4122)
$ws.Range("D17").Value = "Option C"

# Row 18: A18 picks up "Option B" with the same s=5 style.
$ws.Range("B10").Copy()
$ws.Range("A18").PasteSpecial(-4122)
$ws.Range("A18").Value = "Option B"

# Row 22: A22 picks up "Option A" using the "correctStyle" (s=6) already
# used by the numeric "Wrong" column; B22 keeps "Option D".
$ws.Range("C10").Copy()
$ws.Range("A22").PasteSpecial(-4122)
$ws.Range("A22").Value = "Option A"

# Row 29: A29 picks up "Option D" with the same s=5 style as rows 17/18.
$ws.Range("B10").Copy()
$ws.Range("A29").PasteSpecial(-4122)
$ws.Range("A29").Value = "Option D"

# Row 19-21: drop column D:E entirely (no replacement value).
$ws.Range("D19:E21").Clear()

# Row 22: drop column D:E (A22/B22 handled above).
$ws.Range("D22:E22").Clear()

# Rows 23-40: drop column D:E entirely.
$ws.Range("D23:E28").Clear()
$ws.Range("D29:E29").Clear()
$ws.Range("D30:E40").Clear()

Write-Host "done"
